$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: nombre_aides 153 -> 154, montant_total 601891.72 -> 607791.72
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = "154"
$c.Style = "Normal"
$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "607791.72"
$d.Style = "Normal"

# Row 36: nombre_aides 757 -> 762, montant_total 3237900.76 -> 3269577.29
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = "762"
$c.Style = "Normal"
$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "3269577.29"
$d.Style = "Normal"

# Row 37: nombre_aides 357 -> 358, montant_total 2660305.18 -> 2662305.18
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "358"
$c.Style = "Normal"
$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "2662305.18"
$d.Style = "Normal"

# Row 43: nombre_aides 230 -> 233, montant_total 762938.81 -> 806226.81
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "233"
$c.Style = "Normal"
$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "806226.81"
$d.Style = "Normal"

# Row 44: nombre_aides 98 -> 99, montant_total 432429.00 -> 477429.00
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = "99"
$c.Style = "Normal"
$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "477429.00"
$d.Style = "Normal"

# Row 45: nombre_aides 28 -> 31, montant_total 208095.14 -> 308095.14
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = "31"
$c.Style = "Normal"
$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "308095.14"
$d.Style = "Normal"

# Row 47: nombre_aides 93 -> 95, montant_total 292908.00 -> 313908.00
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "95"
$c.Style = "Normal"
$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "313908.00"
$d.Style = "Normal"

# Row 62: nombre_aides 1156 -> 1157, montant_total 3677757.20 -> 3687757.20
$c = $ws.Range("C62")
$c.NumberFormat = "@"
$c.Value = "1157"
$c.Style = "Normal"
$d = $ws.Range("D62")
$d.NumberFormat = "@"
$d.Value = "3687757.20"
$d.Style = "Normal"

# Row 64: nombre_aides 3186 -> 3188, montant_total 19032035.47 -> 19048085.47
$c = $ws.Range("C64")
$c.NumberFormat = "@"
$c.Value = "3188"
$c.Style = "Normal"
$d = $ws.Range("D64")
$d.NumberFormat = "@"
$d.Value = "19048085.47"
$d.Style = "Normal"

# Row 65: nombre_aides 1104 -> 1105, montant_total 8071490.96 -> 8081490.96
$c = $ws.Range("C65")
$c.NumberFormat = "@"
$c.Value = "1105"
$c.Style = "Normal"
$d = $ws.Range("D65")
$d.NumberFormat = "@"
$d.Value = "8081490.96"
$d.Style = "Normal"

# Row 110: nombre_aides 408 -> 409, montant_total 1255073.68 -> 1277943.68
$c = $ws.Range("C110")
$c.NumberFormat = "@"
$c.Value = "409"
$c.Style = "Normal"
$d = $ws.Range("D110")
$d.NumberFormat = "@"
$d.Value = "1277943.68"
$d.Style = "Normal"

# Row 111: nombre_aides 1648 -> 1657, montant_total 6398561.19 -> 6572255.16
$c = $ws.Range("C111")
$c.NumberFormat = "@"
$c.Value = "1657"
$c.Style = "Normal"
$d = $ws.Range("D111")
$d.NumberFormat = "@"
$d.Value = "6572255.16"
$d.Style = "Normal"

# Row 112: nombre_aides 655 -> 661, montant_total 3942575.74 -> 4103463.08
$c = $ws.Range("C112")
$c.NumberFormat = "@"
$c.Value = "661"
$c.Style = "Normal"
$d = $ws.Range("D112")
$d.NumberFormat = "@"
$d.Value = "4103463.08"
$d.Style = "Normal"

# Row 113: nombre_aides 211 -> 215, montant_total 2088709.73 -> 2241621.73
$c = $ws.Range("C113")
$c.NumberFormat = "@"
$c.Value = "215"
$c.Style = "Normal"
$d = $ws.Range("D113")
$d.NumberFormat = "@"
$d.Value = "2241621.73"
$d.Style = "Normal"

# Row 114: nombre_aides 71 -> 74, montant_total 810159.00 -> 931840.00
$c = $ws.Range("C114")
$c.NumberFormat = "@"
$c.Value = "74"
$c.Style = "Normal"
$d = $ws.Range("D114")
$d.NumberFormat = "@"
$d.Value = "931840.00"
$d.Style = "Normal"

# Row 115: nombre_aides 8 -> 10, montant_total 200000.00 -> 254339.00
$c = $ws.Range("C115")
$c.NumberFormat = "@"
$c.Value = "10"
$c.Style = "Normal"
$d = $ws.Range("D115")
$d.NumberFormat = "@"
$d.Value = "254339.00"
$d.Style = "Normal"
